# Example.xlsx update:
#  1. add json value "type" fixed: sheet2!B2 "type" -> "name1", plus new
#     "name2" / "timestamp_arry" / "timestamp[]" / timestamp-array sample
#     columns/rows.
#  2. add test case: a new row (id=2) with the double/"clip" style nested
#     array and a new timestamp array sample.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: rework the small "TestCase" table -----------------------

# Row 2 (field names): B2 used to say "type" - rename it, and extend the
# table with two more named columns.
$ws2.Range("B2").Value = "name1"
$ws2.Range("C2").Value = "name2"

# Row 5: brand new sample row (id=2).
$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "[[1.0,0.1],[1,1],[2,2]]"

# Finish row 2 (new "timestamp array" column header) and row 3 (its type).
$ws2.Range("D2").Value = "timestamp_arry"
$ws2.Range("D3").Value = "timestamp[]"

# Row 4 (sample row, id=1): keep existing B4/C4, add sample timestamp array.
$ws2.Range("D4").Value = "[""2018/01/01 23:59:59""]"

$ws2.Range("D5").Value = "[""2018/12/31 00:00:00""]"

# --- View / active sheet & selection housekeeping ---------------------

# Sheet1 keeps the same selection anchor as before but one row further
# down (W8 -> W9), and is no longer the active tab.
$ws1.Activate()
$ws1.Range("W9").Select()

# Sheet2 becomes the active/selected tab, with its own new selection.
$ws2.Activate()
$ws2.Range("D6").Select()
